$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the two target paragraphs by their (stable, pre-edit) opening text.
# ---------------------------------------------------------------------------
$paraMark = $null
$paraPeriod = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "La dificultad que he encontrado a la hora de realizar este requisito*") {
        $paraMark = $p
    }
    if ($t -like "La dificultad que he encontrado en este requisito es similar*") {
        $paraPeriod = $p
    }
}

# ---------------------------------------------------------------------------
# Section 1: "...atributo <b>mark</b>..." paragraph
# ---------------------------------------------------------------------------

# 1) "...del desarrollo del atributo " -> "...del desarrollo del atributo derivado "
$r1 = $paraMark.Range
$r1.Find.ClearFormatting()
$r1.Find.Execute("del desarrollo del atributo ", $true, $false, $false, $false, $false, $true, 1, $false, `
    "del desarrollo del atributo derivado ", 2) | Out-Null

# 2) the bold run "mark" becomes bold *and* italic
$r1 = $paraMark.Range
$r1.Find.ClearFormatting()
$r1.Find.Text = "mark"
$r1.Find.Forward = $true
$r1.Find.Wrap = 1
$r1.Find.MatchCase = $true
$r1.Find.MatchWholeWord = $true
if ($r1.Find.Execute()) {
    $r1.Italic = $true
}

# 3) remove ", al tratarse también de un atributo derivado" (leaving the full stop).
#    The preceding run "(conglomerado)" shares identical formatting with this run, and
#    this engine coalesces adjacent same-format runs whenever a length-changing edit
#    touches one of them. Shield "(conglomerado)" with a throw-away Bold flag while we
#    shrink the text, then put it back the way it was so the two stay separate runs
#    (matching the reference XML, where they remain distinct <w:r> elements).
$rShield = $paraMark.Range
$rShield.Find.ClearFormatting()
$rShield.Find.Text = "(conglomerado)"
$rShield.Find.Forward = $true
$rShield.Find.Wrap = 0
$rShield.Find.MatchCase = $true
if ($rShield.Find.Execute()) {
    $rShield.Bold = $true
}

$r1 = $paraMark.Range
$r1.Find.ClearFormatting()
$r1.Find.Execute(", al tratarse también de un atributo derivado. No se dieron", $true, $false, $false, $false, $false, $true, 1, $false, `
    ". No se dieron", 2) | Out-Null

$rUnshield = $paraMark.Range
$rUnshield.Find.ClearFormatting()
$rUnshield.Find.Text = "(conglomerado)"
$rUnshield.Find.Forward = $true
$rUnshield.Find.Wrap = 0
$rUnshield.Find.MatchCase = $true
if ($rUnshield.Find.Execute()) {
    $rUnshield.Bold = $false
}

# 4) final sentence tweaks (all inside the last, plain run of the paragraph)
$r1 = $paraMark.Range
$r1.Find.ClearFormatting()
$r1.Find.Execute("solución aportada entendí", $true, $false, $false, $false, $false, $true, 1, $false, `
    "solución aportada por el tutor entendí", 2) | Out-Null

$r1 = $paraMark.Range
$r1.Find.ClearFormatting()
$r1.Find.Execute("se implementa en la capa de servicios y la teoría", $true, $false, $false, $false, $false, $true, 1, $false, `
    "se implementa en la capa de servicios en el método bind/unbind y la teoría", 2) | Out-Null

$r1 = $paraMark.Range
$r1.Find.ClearFormatting()
$r1.Find.Execute("pertenece a la L03.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "pertenece al contenido L03.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Section 2: "...atributo <b><i>period</i></b>..." paragraph
# ---------------------------------------------------------------------------

# 1) Rewrite the opening sentence and drop the "follow up" run entirely
#    (it gets folded into plain text as part of the rewritten sentence).
$r2 = $paraPeriod.Range
$r2.Find.ClearFormatting()
$r2.Find.Execute( `
    "La dificultad que he encontrado en este requisito es similar a la del anterior requisito, pues no se especifica de manera clara y concisa ni en la teoría ni en nuestras preguntas realizadas en el último follow up el cómo calcular la validación del atributo ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "La dificultad que he encontrado en este requisito es similar a la del anterior requisito, aunque esta vez se especifica en la teoría el cómo calcular la validación del atributo ", `
    2) | Out-Null

# 2) rewrite the tail after "period": split the old trailing sentence, add a
#    new one about AuditRecord.java, bold+italic.
$r2 = $paraPeriod.Range
$r2.Find.ClearFormatting()
$r2.Find.Execute( `
    ", donde su duración tiene que ser como mínimo de una hora por lo que entiendo que su implementación dependerá de la futura teoría dada en L03 en la capa de servicio en el método que gestione las validaciones, por lo que dicha validación aún no podrá ser realizada. El mismo enlace del requisito anterior podría responder este problema ya que se trata de un problema similar en la capa de servicio.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ", donde su duración tiene que ser como mínimo de una hora, por lo que entiendo que su implementación dependerá de la futura teoría dada en L03 en la capa de servicio en el método que gestione las validaciones, por lo que dicha validación aún no podrá ser realizada. Aún así se ha dejado comentado un método en la clase @@AUDITRECORD@@ el cual sirve de base para resolver este problema futuro, donde en cualquier caso se puede copiar y pegar donde haga falta haciendo mínimas modificaciones para que se cumpla la validación.", `
    2) | Out-Null

# 3) give "@@AUDITRECORD@@" its final text/formatting (bold + italic)
$r2 = $paraPeriod.Range
$r2.Find.ClearFormatting()
$r2.Find.Execute("@@AUDITRECORD@@", $true, $false, $false, $false, $false, $true, 1, $false, `
    "AuditRecord.java", 2) | Out-Null

$r2 = $paraPeriod.Range
$r2.Find.ClearFormatting()
$r2.Find.Text = "AuditRecord.java"
$r2.Find.Forward = $true
$r2.Find.Wrap = 1
$r2.Find.MatchCase = $true
$r2.Find.MatchWholeWord = $true
if ($r2.Find.Execute()) {
    $r2.Bold = $true
    $r2.Italic = $true
}

Write-Output "done"
